# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" sheet with a newer data snapshot (19:52 instead
# of 19:22) and the corresponding updated per-country figures. A few countries
# (Irlanda, Jordania, Consejo Danes para los Refugiados) overtake their
# neighbour in the case-count ranking, so those row pairs swap contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last-updated timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 19:52"

# Estados Unidos
$ws.Range("B4").Value = 771214
$ws.Range("C4").Value = 6578
$ws.Range("E4").Value = 658369

# Francia
$ws.Range("B7").Value = 155383
$ws.Range("C7").Value = 2489
$ws.Range("D7").Value = 37409
$ws.Range("E7").Value = 97709
$ws.Range("F7").Value = 5683
$ws.Range("G7").Value = 547
$ws.Range("H7").Value = 20265

# Alemania
$ws.Range("B8").Value = 146398
$ws.Range("C8").Value = 656
$ws.Range("E8").Value = 50192
$ws.Range("G8").Value = 64
$ws.Range("H8").Value = 4706

# Turquia
$ws.Range("B10").Value = 90980
$ws.Range("C10").Value = 4674
$ws.Range("D10").Value = 13430
$ws.Range("E10").Value = 75410
$ws.Range("F10").Value = 1909
$ws.Range("G10").Value = 123
$ws.Range("H10").Value = 2140

# Irlanda overtakes Peru -> the two rows swap country + figures
$ws.Range("A21").Value = "Irlanda"
$ws.Range("B21").Value = 15652
$ws.Range("C21").Value = 401
$ws.Range("D21").Value = 77
$ws.Range("E21").Value = 14888
$ws.Range("F21").Value = 294
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = 687

$ws.Range("A22").Value = "Peru"
$ws.Range("B22").Value = 15628
$ws.Range("D22").Value = 6811
$ws.Range("E22").Value = 8417
$ws.Range("F22").Value = 167
$ws.Range("H22").Value = 400

# Singapur
$ws.Range("F35").Value = 23

# Jordania overtakes Taiwan -> the two rows swap country + figures
$ws.Range("A107").Value = "Jordania"
$ws.Range("B107").Value = 425
$ws.Range("C107").Value = 8
$ws.Range("D107").Value = 282
$ws.Range("E107").Value = 136
$ws.Range("F107").Value = 5
$ws.Range("H107").Value = 7

$ws.Range("A108").Value = "Taiwan"
$ws.Range("B108").Value = 422
$ws.Range("C108").Value = 2
$ws.Range("D108").Value = 203
$ws.Range("E108").Value = 213
$ws.Range("F108").Value = 0
$ws.Range("H108").Value = 6

# Consejo Danes para los Refugiados overtakes Mauricio -> rows swap
$ws.Range("A112").Value = "Consejo Danes para los Refugiados"
$ws.Range("B112").Value = 332
$ws.Range("C112").Value = 5
$ws.Range("D112").Value = 27
$ws.Range("E112").Value = 280
$ws.Range("F112").Value = 0
$ws.Range("H112").Value = 25

$ws.Range("A113").Value = "Mauricio"
$ws.Range("B113").Value = 328
$ws.Range("D113").Value = 224
$ws.Range("E113").Value = 95
$ws.Range("F113").Value = 3
$ws.Range("H113").Value = 9
